$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get new values that Excel would otherwise auto-parse as
# numbers (dropping trailing zeros / significant digits), so force them to
# keep a Text format first, preserving the exact string representation.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values (coin name/link swaps and refreshed price/volume figures).
$ws.Range("D2").Value = "43.647.45"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "2.364.96"
$ws.Range("E3").Value = "  +6.21%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "315.07"
$ws.Range("E5").Value = "  +6.09%  "
$ws.Range("D6").Value = "108.40"
$ws.Range("E6").Value = "  -3.07%  "
$ws.Range("D7").Value = "0.641"
$ws.Range("E7").Value = "  +2.44%  "
$ws.Range("D9").Value = "0.639"
$ws.Range("E9").Value = "  +5.21%  "
$ws.Range("D10").Value = "43.05"
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("D11").Value = "0.0940"
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").Value = "8.78"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").Value = "1.04"
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D16").Value = "2.720.24"
$ws.Range("E16").Value = "  +6.25%  "
$ws.Range("D17").Value = "2.374.09"
$ws.Range("E17").Value = "  +6.16%  "
$ws.Range("D18").Value = "43.629.89"
$ws.Range("E18").Value = "  +2.78%  "
$ws.Range("E19").Value = "  +2.83%  "
$ws.Range("D20").Value = "7.28"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").Value = "75.27"
$ws.Range("E21").Value = "  +3.36%  "
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").Value = "2.54"
$ws.Range("E23").Value = "  +7.36%  "
$ws.Range("D24").Value = "257.85"
$ws.Range("E24").Value = "  +12.32%  "
$ws.Range("D25").Value = "9.34"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").Value = "12.06"
$ws.Range("E26").Value = "  +2.69%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "39.16"
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").Value = "22.74"
$ws.Range("E30").Value = "  +7.78%  "
$ws.Range("D31").Value = "3.21"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").Value = "173.30"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "0.0926"
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("D34").Value = "6.00"
$ws.Range("E34").Value = "  +5.35%  "
$ws.Range("D35").Value = "0.132"
$ws.Range("E35").Value = "  +4.94%  "
$ws.Range("D36").Value = "4.96"
$ws.Range("E36").Value = "  -5.26%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0376"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "4.13"
$ws.Range("E38").Value = "  -5.82%  "
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "2.82"
$ws.Range("E40").Value = "  +15.25%  "
$ws.Range("D41").Value = "1.50"
$ws.Range("E41").Value = "  +13.35%  "
$ws.Range("D42").Value = "71.86"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").Value = "0.233"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D44").Value = "12.85"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "5.65"
$ws.Range("E46").Value = "  +3.06%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "9.31"
$ws.Range("E47").Value = "  +9.68%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "111.77"
$ws.Range("E48").Value = "  +8.23%  "
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").Value = "0.100"
$ws.Range("E50").Value = "  +2.71%  "
$ws.Range("D51").Value = "0.479"
$ws.Range("E51").Value = "  +8.21%  "
